$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-09-04 22:49:34"
$wsDeDe.Range("H4").Value = "2016-09-04 22:49:34"

$wsZhCn.Range("H4").Value = "2016-09-04 22:49:29"
$wsZhCn.Range("K4").Value = "2016-09-04 22:49:48"

$wsDeDe.Range("K4").Value = "2016-09-04 22:49:55"
